$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1, 1).Value = "username"
$ws.Cells.Item(1, 2).Value = "password"
$ws.Cells.Item(1, 3).Value = "role"
$ws.Cells.Item(1, 4).Value = "area"
$ws.Cells.Item(1, 5).Value = "approved"
$ws.Cells.Item(1, 6).Value = "celular"
$ws.Cells.Item(1, 7).Value = "correo_personal"

# Row 2 - admin
$ws.Cells.Item(2, 1).Value = "admin@veracruz.gob.mx"
$ws.Cells.Item(2, 2).Value = "admin123"
$ws.Cells.Item(2, 3).Value = "admin"
$ws.Cells.Item(2, 4).Value = "Todas"
$ws.Cells.Item(2, 5).Value = $true
$ws.Cells.Item(2, 6).Value = 22223333
$ws.Cells.Item(2, 7).Value = "admin@gmail.com"

# Row 3 - policia municipal
$ws.Cells.Item(3, 1).Value = "policiampal@veracruzmunicipio.gob.mx"
$ws.Cells.Item(3, 2).Value = "Poli123"
$ws.Cells.Item(3, 3).Value = "enlace"
$ws.Cells.Item(3, 4).Value = "COMANDANCIA DE LA POLÍCIA MUNICIPAL"
$ws.Cells.Item(3, 5).Value = $true
$ws.Cells.Item(3, 6).Value = "nan"

# Row 4 - administracion
$ws.Cells.Item(4, 1).Value = "dlopez@veracruzmunicipio.gob.mx"
$ws.Cells.Item(4, 2).Value = "Admon123"
$ws.Cells.Item(4, 3).Value = "enlace"
$ws.Cells.Item(4, 4).Value = "ADMINISTRACIÓN"
$ws.Cells.Item(4, 5).Value = $true
$ws.Cells.Item(4, 6).Value = "nan"

# Row 5 - alumbrado publico
$ws.Cells.Item(5, 1).Value = "diralumbrado@veracruzmunicipio.gob.mx"
$ws.Cells.Item(5, 2).Value = "Alum123"
$ws.Cells.Item(5, 3).Value = "enlace"
$ws.Cells.Item(5, 4).Value = "ALUMBRADO PÚBLICO"
$ws.Cells.Item(5, 5).Value = $true
$ws.Cells.Item(5, 6).Value = "nan"

# Row 6 - archivo municipal
$ws.Cells.Item(6, 1).Value = "siaveracruz2022@gmail.com"
$ws.Cells.Item(6, 2).Value = "Arch123"
$ws.Cells.Item(6, 3).Value = "enlace"
$ws.Cells.Item(6, 4).Value = "ARCHIVO MUNICIPAL"
$ws.Cells.Item(6, 5).Value = $true
$ws.Cells.Item(6, 6).Value = "nan"

# Row 7 - asuntos legales
$ws.Cells.Item(7, 1).Value = "asuntoslegales@veracruzmunicipio.gob.mx"
$ws.Cells.Item(7, 2).Value = "Legal123"
$ws.Cells.Item(7, 3).Value = "enlace"
$ws.Cells.Item(7, 4).Value = "ASUNTOS LEGALES"
$ws.Cells.Item(7, 5).Value = $true
$ws.Cells.Item(7, 6).Value = "nan"

# Column widths (characters) to match the authored sheet
# (engine snaps ColumnWidth to a 1/6-character pixel grid, same as Excel's
# internal storage granularity, so these land on the closest achievable width)
$ws.Columns.Item(1).ColumnWidth = 31.83
$ws.Columns.Item(2).ColumnWidth = 18.67

# Restore selection to A5, as in the authored file
[void]$ws.Range("A5").Select()
